$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.7268139522358
$ws.Range("C2").Value = 9.909499752752273
$ws.Range("D2").Value = 6.200105789641675
$ws.Range("F2").Value = 50.82474115624333
$ws.Range("G2").Value = 3.746847011814705
$ws.Range("I2").Value = 34.55137006959226
$ws.Range("J2").Value = 11.28123371207078
$ws.Range("K2").Value = 13.34476527065363
$ws.Range("M2").Value = 17.52028713280054
$ws.Range("B3").Value = 11.63025426801422
$ws.Range("C3").Value = 9.802888649705475
$ws.Range("D3").Value = 6.2305438245558
$ws.Range("F3").Value = 50.54717205241624
$ws.Range("G3").Value = 3.749938806359614
$ws.Range("I3").Value = 34.40850612635186
$ws.Range("J3").Value = 11.29199934311104
$ws.Range("K3").Value = 13.28843043936431
$ws.Range("M3").Value = 17.51878735781905
$ws.Range("B4").Value = 11.57483388182553
$ws.Range("C4").Value = 9.739888868531974
$ws.Range("D4").Value = 6.250861242848331
$ws.Range("F4").Value = 50.38344539638445
$ws.Range("G4").Value = 3.751935838986771
$ws.Range("I4").Value = 34.32428382170799
$ws.Range("J4").Value = 11.30005132815327
$ws.Range("K4").Value = 13.25773267807585
$ws.Range("M4").Value = 17.52192639610523
$ws.Range("B5").Value = 11.55324927226296
$ws.Range("C5").Value = 9.714863554025175
$ws.Range("D5").Value = 6.259550371940671
$ws.Range("F5").Value = 50.31844795080115
$ws.Range("G5").Value = 3.752774543372936
$ws.Range("I5").Value = 34.29085570904437
$ws.Range("J5").Value = 11.30369521226023
$ws.Range("K5").Value = 13.24621296884049
$ws.Range("M5").Value = 17.52422770922233
$ws.Range("B6").Value = 11.54972629059011
$ws.Range("C6").Value = 9.710748053830429
$ws.Range("D6").Value = 6.26101793948924
$ws.Range("F6").Value = 50.30776024721271
$ws.Range("G6").Value = 3.752915315951314
$ws.Range("I6").Value = 34.28535932567953
$ws.Range("J6").Value = 11.3043221797579
$ws.Range("K6").Value = 13.24436021490056
$ws.Range("M6").Value = 17.52467159950874
$ws.Range("B7").Value = 11.57453870250024
$ws.Range("C7").Value = 9.739548708960044
$ws.Range("D7").Value = 6.250976768759175
$ws.Range("F7").Value = 50.38256179663367
$ws.Range("G7").Value = 3.75194704912166
$ws.Range("I7").Value = 34.32382936468663
$ws.Range("J7").Value = 11.30009900250112
$ws.Range("K7").Value = 13.25757329758272
$ws.Range("M7").Value = 17.5219532926701
$ws.Range("B8").Value = 11.69273433958391
$ws.Range("C8").Value = 9.872247562184764
$ws.Range("D8").Value = 6.210262958394383
$ws.Range("F8").Value = 50.72766743802439
$ws.Range("G8").Value = 3.747892639694648
$ws.Range("I8").Value = 34.50139345859762
$ws.Range("J8").Value = 11.28464650297192
$ws.Range("K8").Value = 13.32454037052882
$ws.Range("M8").Value = 17.51892828130674
$ws.Range("B9").Value = 11.95391327941783
$ws.Range("C9").Value = 10.15064581834924
$ws.Range("D9").Value = 6.143336359493149
$ws.Range("F9").Value = 51.45585340878589
$ws.Range("G9").Value = 3.740720705043262
$ws.Range("I9").Value = 34.87669227022515
$ws.Range("J9").Value = 11.26578057410915
$ws.Range("K9").Value = 13.48618604742958
$ws.Range("M9").Value = 17.54513293732375
$ws.Range("B10").Value = 12.16187044770434
$ws.Range("C10").Value = 10.36439599058239
$ws.Range("D10").Value = 6.102027570888147
$ws.Range("F10").Value = 52.01973410704878
$ws.Range("G10").Value = 3.735920548279336
$ws.Range("I10").Value = 35.16801967089869
$ws.Range("J10").Value = 11.25888496680118
$ws.Range("K10").Value = 13.6225937016153
$ws.Range("M10").Value = 17.5838403845809
$ws.Range("B11").Value = 12.2595419287212
$ws.Range("C11").Value = 10.46322536962601
$ws.Range("D11").Value = 6.084940071573271
$ws.Range("F11").Value = 52.28196053137268
$ws.Range("G11").Value = 3.733837473165198
$ws.Range("I11").Value = 35.30372277087739
$ws.Range("J11").Value = 11.25725790323715
$ws.Range("K11").Value = 13.68827605884849
$ws.Range("M11").Value = 17.6056325095603
$ws.Range("B12").Value = 12.29693235967632
$ws.Range("C12").Value = 10.50084380353529
$ws.Range("D12").Value = 6.078714343019469
$ws.Range("F12").Value = 52.38202957931384
$ws.Range("G12").Value = 3.733063031741848
$ws.Range("I12").Value = 35.35554699956281
$ws.Range("J12").Value = 11.25685853824882
$ws.Range("K12").Value = 13.71365131535275
$ws.Range("M12").Value = 17.61448172996833
$ws.Range("B13").Value = 12.28886227554701
$ws.Range("C13").Value = 10.49273393696941
$ws.Range("D13").Value = 6.08004427328719
$ws.Range("F13").Value = 52.36044453680652
$ws.Range("G13").Value = 3.73322918360244
$ws.Range("I13").Value = 35.3443666484302
$ws.Range("J13").Value = 11.25693491383536
$ws.Range("K13").Value = 13.70816425468289
$ws.Range("M13").Value = 17.61254941966722
$ws.Range("B14").Value = 12.26261017307966
$ws.Range("C14").Value = 10.46631657632485
$ws.Range("D14").Value = 6.084422968906646
$ws.Range("F14").Value = 52.29017815211535
$ws.Range("G14").Value = 3.733773471806644
$ws.Range("I14").Value = 35.30797775473803
$ws.Range("J14").Value = 11.25722070502538
$ws.Range("K14").Value = 13.69035373351242
$ws.Range("M14").Value = 17.6063486038109
$ws.Range("B15").Value = 12.24658154807285
$ws.Range("C15").Value = 10.45015938683919
$ws.Range("D15").Value = 6.087136940968232
$ws.Range("F15").Value = 52.24723661434337
$ws.Range("G15").Value = 3.734108733480042
$ws.Range("I15").Value = 35.28574472098253
$ws.Range("J15").Value = 11.25742397899011
$ws.Range("K15").Value = 13.67950916477618
$ws.Range("M15").Value = 17.60262802650287
$ws.Range("B16").Value = 12.15554601103503
$ws.Range("C16").Value = 10.35796630763528
$ws.Range("D16").Value = 6.103178560627223
$ws.Range("F16").Value = 52.0027075589628
$ws.Range("G16").Value = 3.736058698102728
$ws.Range("I16").Value = 35.15921335290916
$ws.Range("J16").Value = 11.2590216526361
$ws.Range("K16").Value = 13.61837262856563
$ws.Range("M16").Value = 17.58250002738321
$ws.Range("B17").Value = 12.10045887614418
$ws.Range("C17").Value = 10.30179185910442
$ws.Range("D17").Value = 6.113455956981382
$ws.Range("F17").Value = 51.85412387509862
$ws.Range("G17").Value = 3.737280629342404
$ws.Range("I17").Value = 35.0823898445996
$ws.Range("J17").Value = 11.26038824624236
$ws.Range("K17").Value = 13.58178369359662
$ws.Range("M17").Value = 17.57122065360262
$ws.Range("B18").Value = 12.06906551942683
$ws.Range("C18").Value = 10.26963388889563
$ws.Range("D18").Value = 6.119527646466423
$ws.Range("F18").Value = 51.769203546152
$ws.Range("G18").Value = 3.737992920173671
$ws.Range("I18").Value = 35.03850328441062
$ws.Range("J18").Value = 11.26131640682217
$ws.Range("K18").Value = 13.56108157819531
$ws.Range("M18").Value = 17.56512718147389
$ws.Range("B19").Value = 12.05848742368636
$ws.Range("C19").Value = 10.25877294799096
$ws.Range("D19").Value = 6.121610968605211
$ws.Range("F19").Value = 51.74054547499022
$ws.Range("G19").Value = 3.738235718355454
$ws.Range("I19").Value = 35.02369623307951
$ws.Range("J19").Value = 11.26165508655604
$ws.Range("K19").Value = 13.55413165469216
$ws.Range("M19").Value = 17.5631318588972
$ws.Range("B20").Value = 12.10629311428009
$ws.Range("C20").Value = 10.30775624010259
$ws.Range("D20").Value = 6.112345310262349
$ws.Range("F20").Value = 51.86988523555558
$ws.Range("G20").Value = 3.737149573352592
$ws.Range("I20").Value = 35.09053690486193
$ws.Range("J20").Value = 11.2602280621424
$ws.Range("K20").Value = 13.58564329818422
$ws.Range("M20").Value = 17.57238060338089
$ws.Range("B21").Value = 12.27031037697351
$ws.Range("C21").Value = 10.4740710135363
$ws.Range("D21").Value = 6.083130193274886
$ws.Range("F21").Value = 52.31079662864576
$ws.Range("G21").Value = 3.733613211727879
$ws.Range("I21").Value = 35.31865437737039
$ws.Range("J21").Value = 11.2571308811831
$ws.Range("K21").Value = 13.69557162793059
$ws.Range("M21").Value = 17.6081537707653
$ws.Range("B22").Value = 12.37984447119269
$ws.Range("C22").Value = 10.5838812744553
$ws.Range("D22").Value = 6.065464044002273
$ws.Range("F22").Value = 52.60341832551224
$ws.Range("G22").Value = 3.731385738649885
$ws.Range("I22").Value = 35.4702760619393
$ws.Range("J22").Value = 11.25636998678453
$ws.Range("K22").Value = 13.77033713432018
$ws.Range("M22").Value = 17.63501089752322
$ws.Range("B23").Value = 12.32118243725941
$ws.Range("C23").Value = 10.52518314519569
$ws.Range("D23").Value = 6.074762212091319
$ws.Range("F23").Value = 52.44685005653991
$ws.Range("G23").Value = 3.732566947922516
$ws.Range("I23").Value = 35.38912781962709
$ws.Range("J23").Value = 11.25666062676195
$ws.Range("K23").Value = 13.73017269036378
$ws.Range("M23").Value = 17.62036025475588
$ws.Range("B24").Value = 12.10365458966642
$ws.Range("C24").Value = 10.30505931258146
$ws.Range("D24").Value = 6.112846925675584
$ws.Range("F24").Value = 51.86275795552889
$ws.Range("G24").Value = 3.737208793275752
$ws.Range("I24").Value = 35.0868527442227
$ws.Range("J24").Value = 11.26030003752287
$ws.Range("K24").Value = 13.58389733115663
$ws.Range("M24").Value = 17.57185497123794
$ws.Range("B25").Value = 11.88029560084503
$ws.Range("C25").Value = 10.07357579969774
$ws.Range("D25").Value = 6.160060339942413
$ws.Range("F25").Value = 51.25360011444642
$ws.Range("G25").Value = 3.742578122086319
$ws.Range("I25").Value = 34.77235459047488
$ws.Range("J25").Value = 11.2696604263768
$ws.Range("K25").Value = 13.43929184296304
$ws.Range("M25").Value = 17.53461513733536
